$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Albahaca (Terminal La Palmera de La Serena)
# is inserted at row 187, pushing the existing rows 187-223 down to 188-224.
$ws.Rows.Item(187).Insert()

$ws.Cells.Item(187, 1).Value = 8
$ws.Cells.Item(187, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(187, 3).Value = "Coquimbo"
$ws.Cells.Item(187, 4).Value = 45244
$ws.Cells.Item(187, 5).Value = 4
$ws.Cells.Item(187, 6).Value = 100112052
$ws.Cells.Item(187, 7).Value = "Albahaca"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 1000
$ws.Cells.Item(187, 11).Value = 3500
$ws.Cells.Item(187, 12).Value = 4000
$ws.Cells.Item(187, 13).Value = 3750
$ws.Cells.Item(187, 14).Value = "$/paquete"
$ws.Cells.Item(187, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(187, 16).Value = 3750
$ws.Cells.Item(187, 17).Value = 1
$ws.Cells.Item(187, 18).Value = "Hortaliza"
